# Fill in the student-identification row of the first table on the cover
# page: student name, control number, and project number.
$d = $word.ActiveDocument

$t = $d.Tables.Item(1)
$t.Cell(2, 1).Range.Text = "Axel Martin Vega Espinoza"
$t.Cell(2, 2).Range.Text = "20120168"
$t.Cell(2, 3).Range.Text = "4"
